$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 47620336
$ws.Range("I9").Value = 333333340
$ws.Range("K9").Value = 333333340
$ws.Range("M9").Value = -333333171
$ws.Range("H13").Value = 1906.6666
$ws.Range("J13").Value = 1906.6666
$ws.Range("L13").Value = 1906.6666
$ws.Range("N13").Value = -2244.6666
$ws.Range("H112").Value = 1530.4482
$ws.Range("J112").Value = 1530.4482
$ws.Range("L112").Value = 4591.3446
$ws.Range("N112").Value = -6807.3446
$ws.Range("H138").Value = 66736210
$ws.Range("I138").Value = 168916.83
$ws.Range("J138").Value = 111114400
$ws.Range("K138").Value = 506750.49
$ws.Range("L138").Value = 333343200
$ws.Range("M138").Value = -501610.49
$ws.Range("N138").Value = -333353480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 299.5
$ws.Range("K5").Value = 299.5
$ws.Range("M5").Value = -187.5
$ws.Range("H46").Value = 16787
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 16787
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 16787
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -17425
$ws.Range("H63").Value = 4422.5386
$ws.Range("I63").Value = 3500
$ws.Range("K63").Value = 3500
$ws.Range("M63").Value = -2814
$ws.Range("H66").Value = 4422.5386
$ws.Range("I66").Value = 3500
$ws.Range("K66").Value = 17500
$ws.Range("M66").Value = -14068

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 299.5
$ws.Range("K4").Value = 299.5
$ws.Range("M4").Value = -184.5
$ws.Range("H82").Value = 22584.076
$ws.Range("I82").Value = 3857.4
$ws.Range("K82").Value = 3857.4
$ws.Range("M82").Value = -3474.4
$ws.Range("H85").Value = 22584.076
$ws.Range("I85").Value = 3857.4
$ws.Range("K85").Value = 3857.4
$ws.Range("M85").Value = -2531.4
$ws.Range("H97").Value = 5552.5
$ws.Range("I97").Value = 5552.5
$ws.Range("K97").Value = 5552.5
$ws.Range("M97").Value = -4561.5
$ws.Range("H102").Value = 10414
$ws.Range("I102").Value = 10414
$ws.Range("K102").Value = 10414
$ws.Range("M102").Value = -7169

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3049.75
$ws.Range("J16").Value = 3233
$ws.Range("L16").Value = 3233
$ws.Range("N16").Value = -3807
$ws.Range("H31").Value = 2100.9512
$ws.Range("I31").Value = 1370.1111
$ws.Range("J31").Value = 2672.913
$ws.Range("K31").Value = 1370.1111
$ws.Range("L31").Value = 2672.913
$ws.Range("M31").Value = -1075.1111
$ws.Range("N31").Value = -3262.913
$ws.Range("H34").Value = 2100.9512
$ws.Range("I34").Value = 1370.1111
$ws.Range("J34").Value = 2672.913
$ws.Range("K34").Value = 1370.1111
$ws.Range("L34").Value = 2672.913
$ws.Range("M34").Value = -1168.1111
$ws.Range("N34").Value = -3076.913
$ws.Range("H107").Value = 1200
$ws.Range("I107").Value = 523.8095
$ws.Range("J107").Value = 4750
$ws.Range("K107").Value = 523.8095
$ws.Range("L107").Value = 4750
$ws.Range("M107").Value = 1396.1905
$ws.Range("N107").Value = -8590
$ws.Range("H113").Value = 3049.75
$ws.Range("J113").Value = 3233
$ws.Range("L113").Value = 3233
$ws.Range("N113").Value = -7573

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 218.5
$ws.Range("I2").Value = 106
$ws.Range("K2").Value = 636
$ws.Range("M2").Value = -523
$ws.Range("H12").Value = 80.1579
$ws.Range("I12").Value = 66.28570999999999
$ws.Range("J12").Value = 88.25
$ws.Range("K12").Value = 198.85713
$ws.Range("L12").Value = 264.75
$ws.Range("M12").Value = -25.85712999999998
$ws.Range("N12").Value = -610.75
$ws.Range("H37").Value = 74999.75
$ws.Range("J37").Value = 74999.75
$ws.Range("L37").Value = 224999.25
$ws.Range("N37").Value = -225223.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1283.4
$ws.Range("I2").Value = 309.2857
$ws.Range("K2").Value = 309.2857
$ws.Range("M2").Value = -196.2857
$ws.Range("H13").Value = 400.33334
$ws.Range("I13").Value = 400.33334
$ws.Range("K13").Value = 400.33334
$ws.Range("M13").Value = -261.33334
$ws.Range("H41").Value = 3000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws.Range("H58").Value = 19028.625
$ws.Range("I58").Value = 21234.75
$ws.Range("J58").Value = 16822.5
$ws.Range("K58").Value = 21234.75
$ws.Range("L58").Value = 16822.5
$ws.Range("M58").Value = -20957.75
$ws.Range("N58").Value = -17376.5
$ws.Range("H80").Value = 3449.1428
$ws.Range("I80").Value = 2079
$ws.Range("J80").Value = 6874.5
$ws.Range("K80").Value = 2079
$ws.Range("L80").Value = 6874.5
$ws.Range("M80").Value = -1081
$ws.Range("N80").Value = -8870.5
$ws.Range("H83").Value = 3449.1428
$ws.Range("I83").Value = 2079
$ws.Range("J83").Value = 6874.5
$ws.Range("K83").Value = 10395
$ws.Range("L83").Value = 34372.5
$ws.Range("M83").Value = -5403
$ws.Range("N83").Value = -44356.5
$ws.Range("H99").Value = 8072.385
$ws.Range("I99").Value = 8072.385
$ws.Range("K99").Value = 8072.385
$ws.Range("M99").Value = -5826.385
$ws.Range("H113").Value = 3630.5454
$ws.Range("I113").Value = 3320.5
$ws.Range("J113").Value = 4002.6
$ws.Range("K113").Value = 3320.5
$ws.Range("L113").Value = 4002.6
$ws.Range("M113").Value = -1150.5
$ws.Range("N113").Value = -8342.6
$ws.Range("H122").Value = 5387.909
$ws.Range("I122").Value = 6835.3335
$ws.Range("J122").Value = 3651
$ws.Range("K122").Value = 20506.0005
$ws.Range("L122").Value = 10953
$ws.Range("M122").Value = -18056.0005
$ws.Range("N122").Value = -15853
$ws.Range("H132").Value = 7304.115
$ws.Range("I132").Value = 6203.1665
$ws.Range("J132").Value = 9781.25
$ws.Range("K132").Value = 18609.4995
$ws.Range("L132").Value = 29343.75
$ws.Range("M132").Value = -16079.4995
$ws.Range("N132").Value = -34403.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 70553.125
$ws.Range("I22").Value = 1764.1428
$ws.Range("J22").Value = 124055.664
$ws.Range("K22").Value = 1764.1428
$ws.Range("L22").Value = 124055.664
$ws.Range("M22").Value = -1469.1428
$ws.Range("N22").Value = -124645.664
$ws.Range("H27").Value = 70553.125
$ws.Range("I27").Value = 1764.1428
$ws.Range("J27").Value = 124055.664
$ws.Range("K27").Value = 1764.1428
$ws.Range("L27").Value = 124055.664
$ws.Range("M27").Value = -1657.1428
$ws.Range("N27").Value = -124269.664
$ws.Range("H40").Value = 2265915.8
$ws.Range("I40").Value = 59502.5
$ws.Range("J40").Value = 7939550
$ws.Range("K40").Value = 59502.5
$ws.Range("L40").Value = 7939550
$ws.Range("M40").Value = -59366.5
$ws.Range("N40").Value = -7939822
$ws.Range("H46").Value = 21775.4
$ws.Range("I46").Value = 21775.4
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 21775.4
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -21587.4
$ws.Range("N46").Value = ""
$ws.Range("H55").Value = 7671.5
$ws.Range("I55").Value = 1788.8889
$ws.Range("J55").Value = 15234.857
$ws.Range("K55").Value = 1788.8889
$ws.Range("L55").Value = 15234.857
$ws.Range("M55").Value = -1615.8889
$ws.Range("N55").Value = -15580.857
$ws.Range("H107").Value = 3853.5
$ws.Range("I107").Value = 3853.5
$ws.Range("K107").Value = 3853.5
$ws.Range("M107").Value = -1933.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3979
$ws.Range("I126").Value = 3507.6
$ws.Range("J126").Value = 4502.778
$ws.Range("K126").Value = 10522.8
$ws.Range("L126").Value = 13508.334
$ws.Range("M126").Value = -8052.799999999999
$ws.Range("N126").Value = -18448.334
